# Updated cryptos list values (prices & 1h volume %) per upstream diff.
# For numeric-looking "Price" strings we force Text formatting before
# assigning so Excel doesn't auto-convert them to numbers (which would
# lose formatting such as trailing zeros, e.g. "83.30" -> 83.3), then
# restore the cell's original ("Normal") style so no stray formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.076.36"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.979.66"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7348"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3386"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8286"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08099"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "1.985.41"
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.592"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "99.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +11.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "267.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").Value = "31.098.18"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.077"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.80%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008247"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.14%  "
$ws.Range("D21").Value = "2.250.26"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.064"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.977"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.356"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1327"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.594"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.371"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.626"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.412"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05301"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.279"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.789"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02004"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.878"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.773"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4637"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8527"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.644"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.588"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.99%  "
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.957"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +38.76%  "
